$wb = $excel.ActiveWorkbook

# --- Sheet1: remove the extra batch of entries (rows 390-426), trim row 389
#     back to just the running total formula, and drop the old grand-total
#     formula that used to live in row 427 (leaving its styled, empty C cell).
$ws1 = $wb.Worksheets.Item("Sheet1")

# Wipe out the 38 rows that were appended after the original last entry.
$ws1.Range("A390:E426").Clear()

# Row 389 used to be a data row (388th entry) - clear its data cells and
# turn E389 into the new grand-total formula for the remaining entries.
$ws1.Range("A389:D389").Clear()
$ws1.Range("E389").Formula = "=SUM(E2:E388)"

# The old total formula in E427 is gone now; C427 stays as an empty,
# formatted placeholder cell.
$ws1.Range("E427").Clear()

# --- Sheet2 (Diesel): just a selection change.
$ws2 = $wb.Worksheets.Item("Diesel")
$ws2.Activate()
$ws2.Range("N16").Select()

# --- Restore Sheet1 as the active sheet/selection shown in the diff.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 379
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("E390").Select()
